$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "CondFormat" worksheet as a copy of "Tricks" (same data,
#    formulas, column widths), placed as the last (4th) sheet.
# ---------------------------------------------------------------------------
$tricks = $wb.Worksheets.Item("Tricks")
$tricks.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "CondFormat"

# Give it a portrait page setup (creates the <pageSetup> element on save).
$newSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 2. Apply (and then remove) a handful of conditional-formatting rules on the
#    "Tricks" sheet: two "duplicate values" (red) rules and five "greater
#    than/cell value" (green fill) rules. Excel keeps the generated dxf
#    (differential format) records in styles.xml even after the rules that
#    used them are deleted, which is exactly the set of orphaned dxfs we
#    need to end up with.
# ---------------------------------------------------------------------------
$rng = $tricks.Range("A1:F10")

$fcRed1 = $rng.FormatConditions.AddUniqueValues()
$fcRed1.DupeUnique = 1
$fcRed1.Font.Color = 393372
$fcRed1.Interior.Color = 13551615

$fcRed2 = $rng.FormatConditions.AddUniqueValues()
$fcRed2.DupeUnique = 1
$fcRed2.Font.Color = 393372
$fcRed2.Interior.Color = 13551615

$fcGreen1 = $rng.FormatConditions.Add(1, 3, "100")
$fcGreen1.Interior.Color = 5287936

$fcGreen2 = $rng.FormatConditions.Add(1, 3, "100")
$fcGreen2.Interior.Color = 5287936

$fcGreen3 = $rng.FormatConditions.Add(1, 3, "100")
$fcGreen3.Interior.Color = 5287936

$fcGreen4 = $rng.FormatConditions.Add(1, 3, "100")
$fcGreen4.Interior.Color = 5287936

$fcGreen5 = $rng.FormatConditions.Add(1, 3, "100")
$fcGreen5.Interior.Color = 5287936

$rng.FormatConditions.Delete()

# ---------------------------------------------------------------------------
# 3. Fix up selections / active sheet so the saved view matches what the
#    author last had on screen.
# ---------------------------------------------------------------------------

# "Tricks" ends up with its whole table selected (Ctrl+A from A2).
$tricks.Activate()
$tricks.Range("A1:F10").Select()

# The new "CondFormat" sheet becomes the active tab, cursor resting on D11.
$newSheet.Activate()
$newSheet.Range("D11").Select()
